$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.29"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.70%"
$ws.Range("E2").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.92%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.091"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.08%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07825"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.41%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.138"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-19.63%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.798"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.13%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.789"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.73%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9205"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.56%"
$ws.Range("E9").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.14%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07727"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.33%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08883"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.47%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03115"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.62%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1002"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.12%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001520"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.96%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005885"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.22%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.449"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.50%"
$ws.Range("E17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.269"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.94%"
$ws.Range("E18").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1329"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.51%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.151"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.65%"
$ws.Range("E21").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.94%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04589"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.36%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001239"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.33%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004474"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.00%"
$ws.Range("E25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001249"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.18%"
$ws.Range("E26").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01761"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.02%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04770"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "6.06%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007097"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.93%"
$ws.Range("E41").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.56%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002169"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.80%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01082"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.87%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006243"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.75%"
$ws.Range("E45").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.06%"
$ws.Range("E46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003560"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-59.29%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8112"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.15%"
$ws.Range("E48").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.06%"
$ws.Range("E49").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.06%"
$ws.Range("E50").ClearFormats()